$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the data range so numeric-looking strings (e.g. "573.16")
# are stored as text, matching the source inline-string cells, not auto-converted
# to numbers by Excel. Cleared again afterwards so no stray number format lingers.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = "62.491.02"
$ws.Range("E2").Value = "  +3.79%  "
$ws.Range("D3").Value = "2.442.67"
$ws.Range("E3").Value = "  +2.62%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "573.16"
$ws.Range("E5").Value = "  +2.13%  "
$ws.Range("D6").Value = "145.56"
$ws.Range("E6").Value = "  +5.05%  "
$ws.Range("E7").Value = "  -0.24%  "
$ws.Range("E8").Value = "  +2.18%  "
$ws.Range("D9").Value = "2.442.59"
$ws.Range("E9").Value = "  +2.68%  "
$ws.Range("E10").Value = "  +4.45%  "
$ws.Range("E11").Value = "  +0.71%  "
$ws.Range("E12").Value = "  +2.41%  "
$ws.Range("E13").Value = "  +3.92%  "
$ws.Range("D14").Value = "27.41"
$ws.Range("E14").Value = "  +6.30%  "
$ws.Range("E15").Value = "  +5.87%  "
$ws.Range("D16").Value = "2.843.59"
$ws.Range("D17").Value = "62.380.82"
$ws.Range("E17").Value = "  +3.89%  "
$ws.Range("D18").Value = "2.431.62"
$ws.Range("E18").Value = "  +1.89%  "
$ws.Range("D19").Value = "7.91"
$ws.Range("E19").Value = "  -2.36%  "
$ws.Range("E20").Value = "  +4.00%  "
$ws.Range("D21").Value = "327.36"
$ws.Range("E21").Value = "  +1.65%  "
$ws.Range("E22").Value = "  +1.57%  "
$ws.Range("D23").Value = "2.04"
$ws.Range("E23").Value = "  +11.00%  "
$ws.Range("D24").Value = "0.996"
$ws.Range("E24").Value = "  -0.54%  "
$ws.Range("D25").Value = "65.55"
$ws.Range("E25").Value = "  +2.28%  "
$ws.Range("D26").Value = "628.07"
$ws.Range("E26").Value = "  +12.36%  "
$ws.Range("D27").Value = "8.46"
$ws.Range("E27").Value = "  +4.10%  "
$ws.Range("D28").Value = "0.0₃0981"
$ws.Range("E28").Value = "  +5.43%  "
$ws.Range("E29").Value = "  +1.36%  "
$ws.Range("D30").Value = "8.18"
$ws.Range("E30").Value = "  +2.02%  "
$ws.Range("D31").Value = "1.41"
$ws.Range("E31").Value = "  +7.45%  "
$ws.Range("B32").Value = "Kaspa"
$ws.Range("C32").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D32").Value = "0.137"
$ws.Range("E32").Value = "  +4.09%  "
$ws.Range("B33").Value = "PancakeSwap"
$ws.Range("C33").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D33").Value = "1.84"
$ws.Range("E33").Value = "  +2.64%  "
$ws.Range("B34").Value = "BabyDogeCoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D34").Value = "0.0₆0394"
$ws.Range("E34").Value = "  +36.96%  "
$ws.Range("E35").Value = "  +2.70%  "
$ws.Range("D36").Value = "1.00"
$ws.Range("E36").Value = "  +0.08%  "
$ws.Range("D37").Value = "4.75"
$ws.Range("E37").Value = "  +4.11%  "
$ws.Range("E38").Value = "  +1.79%  "
$ws.Range("D39").Value = "151.26"
$ws.Range("E39").Value = "  -0.35%  "
$ws.Range("E40").Value = "  +6.75%  "
$ws.Range("E41").Value = "  +2.51%  "
$ws.Range("D42").Value = "2.75"
$ws.Range("E42").Value = "  +12.98%  "
$ws.Range("E43").Value = "  +5.82%  "
$ws.Range("B45").Value = "WhiteBITCoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D45").Value = "14.81"
$ws.Range("E45").Value = "  +26.60%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "144.84"
$ws.Range("E46").Value = "  +2.86%  "
$ws.Range("B47").Value = "Filecoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D47").Value = "3.59"
$ws.Range("E47").Value = "  +1.88%  "
$ws.Range("D48").Value = "20.53"
$ws.Range("E48").Value = "  +6.90%  "
$ws.Range("E49").Value = "  +1.98%  "
$ws.Range("D50").Value = "0.0516"
$ws.Range("E50").Value = "  +2.85%  "
$ws.Range("D51").Value = "0.0920"
$ws.Range("E51").Value = "  +2.30%  "

$dataRange.ClearFormats()
